# Adobe Journey Optimizer IP Warmup Plan workbook edit
# - Adds a new "OOTB Domain Groups" worksheet with provider domain data
# - Restores "Warmup Plan" as the active/selected tab
# - Updates "Custom Domain Group" sheet selection + header cell style

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "OOTB Domain Groups"

# --- 2. Populate the new sheet with the OOTB domain group table ---
$newSheet.Cells.Item(1, 1).Value = "Gmail"
$newSheet.Cells.Item(1, 2).Value = "gmail.com"
$newSheet.Cells.Item(1, 3).Value = "googlemail.com"
$newSheet.Cells.Item(1, 4).Value = "googlemail.co.uk"
$newSheet.Cells.Item(2, 1).Value = "Microsoft"
$newSheet.Cells.Item(2, 2).Value = "live.com"
$newSheet.Cells.Item(2, 3).Value = "msn.com"
$newSheet.Cells.Item(2, 4).Value = "hotmail.ca"
$newSheet.Cells.Item(2, 5).Value = "hotmail.com"
$newSheet.Cells.Item(2, 6).Value = "hotmail.de"
$newSheet.Cells.Item(2, 7).Value = "hotmail.dk"
$newSheet.Cells.Item(2, 8).Value = "hotmail.co.jp"
$newSheet.Cells.Item(2, 9).Value = "hotmail.it"
$newSheet.Cells.Item(2, 10).Value = "hotmail.es"
$newSheet.Cells.Item(2, 11).Value = "hotmail.fr"
$newSheet.Cells.Item(2, 12).Value = "hotmail.co.uk"
$newSheet.Cells.Item(2, 13).Value = "hotmail.co.kr"
$newSheet.Cells.Item(2, 14).Value = "hotmail.com.au"
$newSheet.Cells.Item(2, 15).Value = "hotmail.com.ar"
$newSheet.Cells.Item(2, 16).Value = "hotmail.co.il"
$newSheet.Cells.Item(2, 17).Value = "hotmail.com.br"
$newSheet.Cells.Item(2, 18).Value = "hotmail.com.tr"
$newSheet.Cells.Item(2, 19).Value = "hotmail.co.th"
$newSheet.Cells.Item(2, 20).Value = "hotmail.jp"
$newSheet.Cells.Item(2, 21).Value = "hotmail.se"
$newSheet.Cells.Item(2, 22).Value = "live.at"
$newSheet.Cells.Item(2, 23).Value = "live.be"
$newSheet.Cells.Item(2, 24).Value = "live.ca"
$newSheet.Cells.Item(2, 25).Value = "live.cl"
$newSheet.Cells.Item(2, 26).Value = "live.cn"
$newSheet.Cells.Item(2, 27).Value = "live.co.kr"
$newSheet.Cells.Item(2, 28).Value = "live.co.uk"
$newSheet.Cells.Item(2, 29).Value = "live.com.ar"
$newSheet.Cells.Item(2, 30).Value = "live.com.au"
$newSheet.Cells.Item(2, 31).Value = "live.com.mx"
$newSheet.Cells.Item(2, 32).Value = "live.com.my"
$newSheet.Cells.Item(2, 33).Value = "live.com.sg"
$newSheet.Cells.Item(2, 34).Value = "live.de"
$newSheet.Cells.Item(2, 35).Value = "live.dk"
$newSheet.Cells.Item(2, 36).Value = "live.fr"
$newSheet.Cells.Item(2, 37).Value = "live.hk"
$newSheet.Cells.Item(2, 38).Value = "live.ie"
$newSheet.Cells.Item(2, 39).Value = "live.in"
$newSheet.Cells.Item(2, 40).Value = "live.it"
$newSheet.Cells.Item(2, 41).Value = "live.jp"
$newSheet.Cells.Item(2, 42).Value = "live.nl"
$newSheet.Cells.Item(2, 43).Value = "live.no"
$newSheet.Cells.Item(2, 44).Value = "live.ru"
$newSheet.Cells.Item(2, 45).Value = "live.se"
$newSheet.Cells.Item(2, 46).Value = "outlook.com"
$newSheet.Cells.Item(2, 47).Value = "live.co.uk"
$newSheet.Cells.Item(2, 48).Value = "hotmail.gr"
$newSheet.Cells.Item(2, 49).Value = "windowslive.com"
$newSheet.Cells.Item(2, 50).Value = "xbox.com"
$newSheet.Cells.Item(2, 51).Value = "hotmail.cl"
$newSheet.Cells.Item(2, 52).Value = "live.at"
$newSheet.Cells.Item(2, 53).Value = "live.jp"
$newSheet.Cells.Item(2, 54).Value = "live.ca"
$newSheet.Cells.Item(2, 55).Value = "hotmail.ca"
$newSheet.Cells.Item(2, 56).Value = "hotmail.ch"
$newSheet.Cells.Item(2, 57).Value = "live.fr"
$newSheet.Cells.Item(2, 58).Value = "live.it"
$newSheet.Cells.Item(2, 59).Value = "live.nl"
$newSheet.Cells.Item(2, 60).Value = "outlook.ie"
$newSheet.Cells.Item(2, 61).Value = "outlook.com.br"
$newSheet.Cells.Item(2, 62).Value = "live.com.pt"
$newSheet.Cells.Item(2, 63).Value = "live.be"
$newSheet.Cells.Item(2, 64).Value = "live.co.za"
$newSheet.Cells.Item(2, 65).Value = "mts.net"
$newSheet.Cells.Item(3, 1).Value = "Yahoo"
$newSheet.Cells.Item(3, 2).Value = "yahoo.com"
$newSheet.Cells.Item(3, 3).Value = "rocketmail.com"
$newSheet.Cells.Item(3, 4).Value = "rogers.com"
$newSheet.Cells.Item(3, 5).Value = "sky.com"
$newSheet.Cells.Item(3, 6).Value = "talk21.com"
$newSheet.Cells.Item(3, 7).Value = "y7mail.com"
$newSheet.Cells.Item(3, 8).Value = "yahoo.at"
$newSheet.Cells.Item(3, 9).Value = "yahoo.be"
$newSheet.Cells.Item(3, 10).Value = "yahoo.bg"
$newSheet.Cells.Item(3, 11).Value = "yahoo.ca"
$newSheet.Cells.Item(3, 12).Value = "yahoo.cl"
$newSheet.Cells.Item(3, 13).Value = "yahoo.co.hu"
$newSheet.Cells.Item(3, 14).Value = "yahoo.co.id"
$newSheet.Cells.Item(3, 15).Value = "yahoo.co.il"
$newSheet.Cells.Item(3, 16).Value = "yahoo.co.in"
$newSheet.Cells.Item(3, 17).Value = "yahoo.co.jp"
$newSheet.Cells.Item(3, 18).Value = "yahoo.co.kr"
$newSheet.Cells.Item(3, 19).Value = "yahoo.com.ar"
$newSheet.Cells.Item(3, 20).Value = "yahoo.com.au"
$newSheet.Cells.Item(3, 21).Value = "yahoo.com.biz"
$newSheet.Cells.Item(3, 22).Value = "yahoo.com.br"
$newSheet.Cells.Item(3, 23).Value = "yahoo.com.cn"
$newSheet.Cells.Item(3, 24).Value = "yahoo.com.co"
$newSheet.Cells.Item(3, 25).Value = "yahoo.com.hk"
$newSheet.Cells.Item(3, 26).Value = "yahoo.com.hr"
$newSheet.Cells.Item(3, 27).Value = "yahoo.com.in"
$newSheet.Cells.Item(3, 28).Value = "yahoo.com.jp"
$newSheet.Cells.Item(3, 29).Value = "yahoo.com.kr"
$newSheet.Cells.Item(3, 30).Value = "yahoo.com.mx"
$newSheet.Cells.Item(3, 31).Value = "yahoo.com.my"
$newSheet.Cells.Item(3, 32).Value = "yahoo.com.net"
$newSheet.Cells.Item(3, 33).Value = "yahoo.com.pe"
$newSheet.Cells.Item(3, 34).Value = "yahoo.com.ph"
$newSheet.Cells.Item(3, 35).Value = "yahoo.com.sg"
$newSheet.Cells.Item(3, 36).Value = "yahoo.com.tr"
$newSheet.Cells.Item(3, 37).Value = "yahoo.com.tw"
$newSheet.Cells.Item(3, 38).Value = "yahoo.com.ua"
$newSheet.Cells.Item(3, 39).Value = "yahoo.com.ve"
$newSheet.Cells.Item(3, 40).Value = "yahoo.com.vn"
$newSheet.Cells.Item(3, 41).Value = "yahoo.co.nz"
$newSheet.Cells.Item(3, 42).Value = "yahoo.co.th"
$newSheet.Cells.Item(3, 43).Value = "yahoo.co.uk"
$newSheet.Cells.Item(3, 44).Value = "yahoo.co.za"
$newSheet.Cells.Item(3, 45).Value = "yahoo.cz"
$newSheet.Cells.Item(3, 46).Value = "yahoo.de"
$newSheet.Cells.Item(3, 47).Value = "yahoo.dk"
$newSheet.Cells.Item(3, 48).Value = "yahoo.ee"
$newSheet.Cells.Item(3, 49).Value = "yahoo.es"
$newSheet.Cells.Item(3, 50).Value = "yahoo.fi"
$newSheet.Cells.Item(3, 51).Value = "yahoo.fr"
$newSheet.Cells.Item(3, 52).Value = "yahoogroups.co.kr"
$newSheet.Cells.Item(3, 53).Value = "yahoogroups.com.cn"
$newSheet.Cells.Item(3, 54).Value = "yahoogroups.com.sg"
$newSheet.Cells.Item(3, 55).Value = "yahoogroups.com.tw"
$newSheet.Cells.Item(3, 56).Value = "yahoogrupper.dk"
$newSheet.Cells.Item(3, 57).Value = "yahoogruppi.it"
$newSheet.Cells.Item(3, 58).Value = "yahoo.gr"
$newSheet.Cells.Item(3, 59).Value = "yahoo.hr"
$newSheet.Cells.Item(3, 60).Value = "yahoo.hu"
$newSheet.Cells.Item(3, 61).Value = "yahoo.ie"
$newSheet.Cells.Item(3, 62).Value = "yahoo.in"
$newSheet.Cells.Item(3, 63).Value = "yahoo.it"
$newSheet.Cells.Item(3, 64).Value = "yahoo.lt"
$newSheet.Cells.Item(3, 65).Value = "yahoo.lv"
$newSheet.Cells.Item(3, 66).Value = "yahoo.nl"
$newSheet.Cells.Item(3, 67).Value = "yahoo.no"
$newSheet.Cells.Item(3, 68).Value = "yahoo.pl"
$newSheet.Cells.Item(3, 69).Value = "yahoo.pt"
$newSheet.Cells.Item(3, 70).Value = "yahoo.ro"
$newSheet.Cells.Item(3, 71).Value = "yahoo.rs"
$newSheet.Cells.Item(3, 72).Value = "yahoo.se"
$newSheet.Cells.Item(3, 73).Value = "yahoo.si"
$newSheet.Cells.Item(3, 74).Value = "yahoo.sk"
$newSheet.Cells.Item(3, 75).Value = "yahooxtra.co.nz"
$newSheet.Cells.Item(3, 76).Value = "ymail.com"
$newSheet.Cells.Item(3, 77).Value = "aol.com"
$newSheet.Cells.Item(3, 78).Value = "aim.com"
$newSheet.Cells.Item(3, 79).Value = "compuserve.com"
$newSheet.Cells.Item(3, 80).Value = "cs.com"
$newSheet.Cells.Item(3, 81).Value = "netscape.com"
$newSheet.Cells.Item(3, 82).Value = "netscape.net"
$newSheet.Cells.Item(3, 83).Value = "wmconnect.com"
$newSheet.Cells.Item(3, 84).Value = "aol.co.uk"
$newSheet.Cells.Item(3, 85).Value = "aol.in"
$newSheet.Cells.Item(3, 86).Value = "aol.de"
$newSheet.Cells.Item(3, 87).Value = "aol.fr"
$newSheet.Cells.Item(3, 88).Value = "aol.nl"
$newSheet.Cells.Item(3, 89).Value = "aol.pl"
$newSheet.Cells.Item(3, 90).Value = "aol.jp"
$newSheet.Cells.Item(3, 91).Value = "aol.es"
$newSheet.Cells.Item(3, 92).Value = "aol.it"
$newSheet.Cells.Item(3, 93).Value = "aol.com.ar"
$newSheet.Cells.Item(3, 94).Value = "aol.fi"
$newSheet.Cells.Item(3, 95).Value = "aol.cl"
$newSheet.Cells.Item(3, 96).Value = "aol.com.co"
$newSheet.Cells.Item(3, 97).Value = "aol.com.ve"
$newSheet.Cells.Item(3, 98).Value = "aol.com.au"
$newSheet.Cells.Item(3, 99).Value = "aol.at"
$newSheet.Cells.Item(3, 100).Value = "aol.be"
$newSheet.Cells.Item(3, 101).Value = "aol.com.br"
$newSheet.Cells.Item(3, 102).Value = "aol.cz"
$newSheet.Cells.Item(3, 103).Value = "aol.dk"
$newSheet.Cells.Item(3, 104).Value = "myaol.jp"
$newSheet.Cells.Item(3, 105).Value = "aolnorge.no"
$newSheet.Cells.Item(3, 106).Value = "aolpolska.pl"
$newSheet.Cells.Item(3, 107).Value = "aolpolcka.pl"
$newSheet.Cells.Item(3, 108).Value = "aolpoland.pl"
$newSheet.Cells.Item(3, 109).Value = "aol.ru"
$newSheet.Cells.Item(3, 110).Value = "aol.kr"
$newSheet.Cells.Item(3, 111).Value = "aol.se"
$newSheet.Cells.Item(3, 112).Value = "aol.ch"
$newSheet.Cells.Item(3, 113).Value = "aol.com.tr"
$newSheet.Cells.Item(3, 114).Value = "aol.co.nz"
$newSheet.Cells.Item(3, 115).Value = "aolchina.com"
$newSheet.Cells.Item(3, 116).Value = "aol.hk"
$newSheet.Cells.Item(3, 117).Value = "aol.tw"
$newSheet.Cells.Item(3, 118).Value = "luckymail.com"
$newSheet.Cells.Item(3, 119).Value = "verizon.net"
$newSheet.Cells.Item(3, 120).Value = "aol.com.mx"
$newSheet.Cells.Item(3, 121).Value = "bellatlantic.net"
$newSheet.Cells.Item(3, 122).Value = "citlink.net"
$newSheet.Cells.Item(3, 123).Value = "frontier.com"
$newSheet.Cells.Item(3, 124).Value = "frontiernet.net"
$newSheet.Cells.Item(3, 125).Value = "games.com"
$newSheet.Cells.Item(3, 126).Value = "goowy.com"
$newSheet.Cells.Item(3, 127).Value = "gte.net"
$newSheet.Cells.Item(3, 128).Value = "love.com"
$newSheet.Cells.Item(3, 129).Value = "verizon.net.in"
$newSheet.Cells.Item(3, 130).Value = "wild4music.com"
$newSheet.Cells.Item(3, 131).Value = "wow.com"
$newSheet.Cells.Item(3, 132).Value = "yahoo.cn"
$newSheet.Cells.Item(3, 133).Value = "yahoo.ne.jp"
$newSheet.Cells.Item(3, 134).Value = "yahoogroups.ca"
$newSheet.Cells.Item(3, 135).Value = "yahoogroups.co.in"
$newSheet.Cells.Item(3, 136).Value = "yahoogroups.co.uk"
$newSheet.Cells.Item(3, 137).Value = "yahoogroups.com"
$newSheet.Cells.Item(3, 138).Value = "yahoogroups.com.au"
$newSheet.Cells.Item(3, 139).Value = "yahoogroups.com.hk"
$newSheet.Cells.Item(3, 140).Value = "yahoogroups.de"
$newSheet.Cells.Item(3, 141).Value = "ybb.ne.jp"
$newSheet.Cells.Item(3, 142).Value = "ygm.com"
$newSheet.Cells.Item(4, 1).Value = "Apple"
$newSheet.Cells.Item(4, 2).Value = "mac.com"
$newSheet.Cells.Item(4, 3).Value = "icloud.com"
$newSheet.Cells.Item(4, 4).Value = "apple.com"
$newSheet.Cells.Item(4, 5).Value = "me.com"
$newSheet.Cells.Item(5, 1).Value = "Comcast"
$newSheet.Cells.Item(5, 2).Value = "comcast.net"
$newSheet.Cells.Item(6, 1).Value = "Orange"
$newSheet.Cells.Item(6, 2).Value = "orange.fr"
$newSheet.Cells.Item(6, 3).Value = "orange.com"
$newSheet.Cells.Item(6, 4).Value = "wanadoo.fr"
$newSheet.Cells.Item(6, 5).Value = "francetelecom.com"
$newSheet.Cells.Item(6, 6).Value = "voila.fr"
$newSheet.Cells.Item(6, 7).Value = "voila.com"
$newSheet.Cells.Item(7, 1).Value = "La Poste"
$newSheet.Cells.Item(7, 2).Value = "laposte.net"
$newSheet.Cells.Item(8, 1).Value = "Italia Online"
$newSheet.Cells.Item(8, 2).Value = "libero.it"
$newSheet.Cells.Item(8, 3).Value = "inwind.it"
$newSheet.Cells.Item(8, 4).Value = "iol.it"
$newSheet.Cells.Item(8, 5).Value = "blu.it"
$newSheet.Cells.Item(8, 6).Value = "giallo.it"
$newSheet.Cells.Item(8, 7).Value = "virgilio.it"
$newSheet.Cells.Item(9, 1).Value = "WP"
$newSheet.Cells.Item(9, 2).Value = "wp.pl"
$newSheet.Cells.Item(9, 3).Value = "o2.pl"
$newSheet.Cells.Item(10, 1).Value = "United Internet"
$newSheet.Cells.Item(10, 2).Value = "web.de"
$newSheet.Cells.Item(10, 3).Value = "gmx.de"
$newSheet.Cells.Item(10, 4).Value = "gmx.ch"
$newSheet.Cells.Item(10, 5).Value = "gmx.net"
$newSheet.Cells.Item(10, 6).Value = "gmx.com"
$newSheet.Cells.Item(10, 7).Value = "gmx.at"
$newSheet.Cells.Item(10, 8).Value = "gmx.fr"
$newSheet.Cells.Item(10, 9).Value = "mail.com"
$newSheet.Cells.Item(10, 10).Value = "1and1.com"
$newSheet.Cells.Item(10, 11).Value = "1und1.de"
$newSheet.Cells.Item(11, 1).Value = "Bigpond"
$newSheet.Cells.Item(11, 2).Value = "bigpond.com"
$newSheet.Cells.Item(11, 3).Value = "bigpond.net.au"
$newSheet.Cells.Item(11, 4).Value = "bigpond.com.au"
$newSheet.Cells.Item(11, 5).Value = "telstra.com"
$newSheet.Cells.Item(11, 6).Value = "bigpond.net"
$newSheet.Cells.Item(12, 1).Value = "Docomo"
$newSheet.Cells.Item(12, 2).Value = "docomo.ne.jp"
$newSheet.Cells.Item(13, 1).Value = "Softbank"
$newSheet.Cells.Item(13, 2).Value = "softbank.ne.jp"
$newSheet.Cells.Item(13, 3).Value = "c.vodafone.ne.jp"
$newSheet.Cells.Item(13, 4).Value = "d.vodafone.ne.jp"
$newSheet.Cells.Item(13, 5).Value = "h.vodafone.ne.jp"
$newSheet.Cells.Item(13, 6).Value = "k.vodafone.ne.jp"
$newSheet.Cells.Item(13, 7).Value = "n.vodafone.ne.jp"
$newSheet.Cells.Item(13, 8).Value = "q.vodafone.ne.jp"
$newSheet.Cells.Item(13, 9).Value = "r.vodafone.ne.jp"
$newSheet.Cells.Item(13, 10).Value = "s.vodafone.ne.jp"
$newSheet.Cells.Item(13, 11).Value = "t.vodafone.ne.jp"
$newSheet.Cells.Item(13, 12).Value = "jp-c.ne.jp"
$newSheet.Cells.Item(13, 13).Value = "jp-d.ne.jp"
$newSheet.Cells.Item(13, 14).Value = "jp-h.ne.jp"
$newSheet.Cells.Item(13, 15).Value = "jp-k.ne.jp"
$newSheet.Cells.Item(13, 16).Value = "jp-n.ne.jp"
$newSheet.Cells.Item(13, 17).Value = "jp-q.ne.jp"
$newSheet.Cells.Item(13, 18).Value = "jp-r.ne.jp"
$newSheet.Cells.Item(13, 19).Value = "jp-s.ne.jp"
$newSheet.Cells.Item(13, 20).Value = "jp-t.ne.jp"
$newSheet.Cells.Item(14, 1).Value = "KDDI"
$newSheet.Cells.Item(14, 2).Value = "au.com"
$newSheet.Cells.Item(14, 3).Value = "ezweb.ne.jp"
$newSheet.Cells.Item(14, 4).Value = "uqmobile.jp"

# Column A autosizing (approximate "best fit" width used by the source file)
$newSheet.Columns.Item(1).ColumnWidth = 12.8

# Leave the new sheet's selection on E8, matching the authored view state
$newSheet.Activate()
$newSheet.Range("E8").Select()

# --- 3. Fix up the "Custom Domain Group" sheet ---
$ws2 = $wb.Worksheets.Item("Custom Domain Group")
$ws2.Activate()

# A1 used a "code" style (Menlo font); switch it to the regular bold header style
$a1 = $ws2.Range("A1")
$a1.Font.Name = "Calibri"
$a1.Font.Color = 0

# Move the active selection to C6
$ws2.Range("C6").Select()

# --- 4. Re-activate "Warmup Plan" so it is the selected tab on open ---
$ws1 = $wb.Worksheets.Item("Warmup Plan")
$ws1.Activate()
